# Adds a new "Count" field/column (column I) to the Shop sheet, mirroring
# the existing NFExport-style metadata rows (Type/Public/Private/Save/
# Cache/Ref/Upload/Desc) and filling every data row with Count = 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column is I (9)
$col = 9

# Copy formatting from column H (8) into the new column I for every row
# that currently holds data, so the new column visually matches the rest
# of the table. Row 9 ("Desc") is intentionally skipped - it stays an
# 8-column row with no Count description.
$srcRange1 = $ws.Range($ws.Cells.Item(1, 8), $ws.Cells.Item(8, 8))
$dstRange1 = $ws.Range($ws.Cells.Item(1, $col), $ws.Cells.Item(8, $col))
[void]$srcRange1.Copy()
[void]$dstRange1.PasteSpecial(-4122)

$srcRange2 = $ws.Range($ws.Cells.Item(10, 8), $ws.Cells.Item(71, 8))
$dstRange2 = $ws.Range($ws.Cells.Item(10, $col), $ws.Cells.Item(71, $col))
[void]$srcRange2.Copy()
[void]$dstRange2.PasteSpecial(-4122)

# Header row
$ws.Cells.Item(1, $col).Value = "Count"
$ws.Cells.Item(1, $col).HorizontalAlignment = 1

# Type row
$ws.Cells.Item(2, $col).Value = "int"

# Public / Private / Save / Cache / Ref / Upload metadata rows
$ws.Cells.Item(3, $col).Value = $false
$ws.Cells.Item(4, $col).Value = $false
$ws.Cells.Item(5, $col).Value = $true
$ws.Cells.Item(6, $col).Value = $false
$ws.Cells.Item(7, $col).Value = $false
$ws.Cells.Item(8, $col).Value = $false

# Desc row (row 9) stays empty - no Count description was supplied.

# Data rows: every shop entry gets Count = 1
for ($r = 10; $r -le 71; $r++) {
    $ws.Cells.Item($r, $col).Value = 1
}

# Leave the cursor highlighting the newly filled Count data range, as the
# original author did after populating the column.
[void]$ws.Range("I10:I71").Select()
